# Add the new "2022-Q3" sheet (quarterly fund-holding breakdown) right after
# the "总计" (Total) summary sheet, and record its summary numbers in the
# "总计" sheet as a new row. Every other quarter sheet keeps its own data and
# simply slides one position to the right to make room.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- 1. Insert the brand-new quarter worksheet right after "总计" ---------
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Header row - identical wording/order used by the other quarterly sheets.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").NumberFormat = "@"

# Match the bold/bordered/centered header style used on the other sheets.
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Single holding row for 2022-Q3.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "001167"
$newSheet.Range("C2").NumberFormat = "@"
$newSheet.Range("C2").Value = "金鹰科技创新股票"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.66"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "94.84"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "5.16"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.1373"
$newSheet.Range("H2").Value = 8

# Match the bold/bordered/centered style used for column A on other sheets.
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# --- 2. Insert a matching summary row into "总计" ---------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Restore column-A styling (bold/bordered/centered) that Insert() doesn't
# carry down automatically to the freshly inserted row.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.14

# Column A is a plain 0-based rank (row 2 -> 0, row 3 -> 1, ...). Re-stamp it
# top to bottom now that one extra row exists.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6
$totalSheet.Range("A9").Value = 7
